# Regenerate orders with updated distance/size codes.
# Distances: D80->D86, D51->D55, D64->D69
# Sizes:     S30->S31
# These substrings appear throughout Condition, Filename_Left,
# Filename_Right, Distance and Size columns (and nowhere else), so a
# straightforward Find&Replace across the used cells reproduces the
# shared-string-table rewrite shown in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Replace("D80", "D86")
$ws.Cells.Replace("D51", "D55")
$ws.Cells.Replace("D64", "D69")
$ws.Cells.Replace("S30", "S31")
